$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New column F: per-patient "score" values for QS 42..48 (rows 34-38), with
# a thick boxed border, centered/wrapped "B Nazanin" text, and an AVERAGE
# summary formula in F43.
# ---------------------------------------------------------------------------

# F43 gets the plain custom number format (no font/border override) - build
# this first so it claims the earliest new cellXfs slot.
$ws.Range("F43").NumberFormat = "[$-3000401]0.##"
$ws.Range("F43").Formula = "=AVERAGE(F34:F38)"

# F35:F38 share the same look: custom font, centered + wrapped text, and a
# medium border on left/right/bottom only (no top - the row above supplies
# the dividing line).
$rngInner = $ws.Range("F35:F38")
$rngInner.Value = @(@(0.62), @(0.73), @(0.83), @(0.82))
$rngInner.NumberFormat = "[$-3000401]0.##"
$rngInner.Font.Size = 11
$rngInner.Font.Name = "B Nazanin"
$rngInner.WrapText = $true
$rngInner.HorizontalAlignment = -4108
$rngInner.VerticalAlignment = -4108
$rngInner.Borders.Weight = -4138
$rngInner.Borders.Item(8).LineStyle = 0

# F34 is the top of the box: same font/format/alignment, but with a full
# medium border on all four sides.
$ws.Range("F34").Value = 0.73
$ws.Range("F34").NumberFormat = "[$-3000401]0.##"
$ws.Range("F34").Font.Size = 11
$ws.Range("F34").Font.Name = "B Nazanin"
$ws.Range("F34").WrapText = $true
$ws.Range("F34").HorizontalAlignment = -4108
$ws.Range("F34").VerticalAlignment = -4108
$ws.Range("F34").Borders.Weight = -4138

# ---------------------------------------------------------------------------
# Row heights for the boxed rows.
# ---------------------------------------------------------------------------
$ws.Rows.Item(33).RowHeight = 17
$ws.Rows.Item(34).RowHeight = 18
$ws.Rows.Item(35).RowHeight = 18
$ws.Rows.Item(36).RowHeight = 18
$ws.Rows.Item(37).RowHeight = 18
$ws.Rows.Item(38).RowHeight = 18

# ---------------------------------------------------------------------------
# Sheet view: scroll/zoom/selection as left by the author.
# ---------------------------------------------------------------------------
$ws.Range("D27").Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 110
